$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"
$zhError = "Handback file name: ghr1xmx4.xz5 is different with handoff file name: 3f3cea27-d577-4679-b966-978b356dd5ef.4638557a3c50524f2663c487b899f9f7238076cd.zh-cn."
$deError = "Handback file name: ghr1xmx4.xz5 is different with handoff file name: 3f3cea27-d577-4679-b966-978b356dd5ef.4638557a3c50524f2663c487b899f9f7238076cd.de-de."

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Note: the ColumnWidth property is specified in "character" units and the
# runtime snaps it to the nearest pixel (using a 6px max-digit-width + 5px
# padding model) when it re-serializes the stored OOXML column width, i.e.
# stored = Round(input * 6 + 5) / 6. To land exactly on a stored width of
# 40 (as in the target file) we request 39.16666667 (= 235/6), which is
# the middle of the input range that rounds to a stored width of 40.
$targetColumnWidth = 39.16666667

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("P3").Value = $zhError
$wsZh.Range("P1").EntireColumn.ColumnWidth = $targetColumnWidth

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("P3").Value = $deError
$wsDe.Range("P1").EntireColumn.ColumnWidth = $targetColumnWidth
